# "resultados atualizados" -- update the figures in the medicamentos
# (medications) summary table: sample size corrected from 365 to 364,
# the "Profilaxia (%)" row re-labelled to "Numero.Medicamentos (%)",
# and several percentages recomputed against the new n.

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

function Update-Cell($row, $col, $old, $new) {
    $cell = $tbl.Cell($row, $col)
    # Cell.Range.Text carries a trailing end-of-cell mark (CR + cell-mark
    # char) after the visible text, so compare only the leading substring.
    $current = $cell.Range.Text
    if ($current.Length -lt $old.Length -or $current.Substring(0, $old.Length) -ne $old) {
        throw "Cell ($row,$col): expected '$old' but found '$current'"
    }
    $cell.Range.Text = $new
}

# n: 365 -> 364
Update-Cell 2 2 "365" "364"

# Row label: Profilaxia (%) -> Numero.Medicamentos (%)
Update-Cell 3 1 "Profilaxia (%)" "Numero.Medicamentos (%)"

# "2" row value: 225 (61.6) -> 224 (61.5)
Update-Cell 6 2 "225 (61.6)" "224 (61.5)"

# Enoxaparina (%): 350 (95.9) -> 349 (95.9)
Update-Cell 10 2 "350 (95.9)" "349 (95.9)"

# Rivaroxabana (%): 81 (22.2) -> 81 (22.3)
Update-Cell 11 2 "81 (22.2)" "81 (22.3)"

# Warfarina (%): 209 (59.2) -> 208 (59.1)
Update-Cell 12 2 "209 (59.2)" "208 (59.1)"
